$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date in column C for all data rows (2-499)
# from serial date 45177 (2023-09-08) to 45178 (2023-09-09).
$ws.Range("C2:C499").Value = 45178
